$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = 8.9765365930052372
$ws.Range("C5").Value = 0.12547388901598394
$ws.Range("D5").Value = 0.54115727430408356

$ws.Range("B7").Value = 133.56744230503554
$ws.Range("C7").Value = 1.7479545256236488
$ws.Range("D7").Value = 15.127369969283384

$ws.Range("B8").Value = 16.930480369719927
$ws.Range("C8").Value = 0.21705873415793656
$ws.Range("D8").Value = 1.9497855105356285

$ws.Range("B11").Value = 26.553493325286336
$ws.Range("C11").Value = 0.40966609969787399
$ws.Range("D11").Value = 2.886306943695141

$ws.Range("B12").Value = 83.84057987725366
$ws.Range("C12").Value = 0.97500462897211038
$ws.Range("D12").Value = 9.8892494991450945

$ws.Range("B14").Value = 52.567424193058606
$ws.Range("C14").Value = 0.51464515846117276
$ws.Range("D14").Value = 4.330391913346479

$ws.Range("B17").Value = 11.151455677570697
$ws.Range("C17").Value = 0.10308033557842679
$ws.Range("D17").Value = 2.1456047502990248

$ws.Range("B18").Value = 16.569595086914074
$ws.Range("C18").Value = 0.1160539407201957
$ws.Range("D18").Value = 5.7085100282054677

$ws.Range("B19").Value = 19.674251497820549
$ws.Range("C19").Value = 0.18806956991621085
$ws.Range("D19").Value = 5.6203248931339909

$ws.Range("B20").Value = 192.90383952728445
$ws.Range("C20").Value = 2.3158259168389383
$ws.Range("D20").Value = 31.561503970055707

$ws.Range("B21").Value = 37.318171949549381
$ws.Range("C21").Value = 0.37469751323425876
$ws.Range("D21").Value = 1.9565750417516437

$ws.Range("B22").Value = 14.744920967796926
$ws.Range("C22").Value = 0.13086494193995291
$ws.Range("D22").Value = 1.825537126535006

$ws.Range("B23").Value = 57.415276521871924
$ws.Range("C23").Value = 0.8358872055547717
$ws.Range("D23").Value = 9.3834857689038067

$ws.Range("B24").Value = 10.71494936093238
$ws.Range("C24").Value = 0.16172053508639397
$ws.Range("D24").Value = 0.81043870275297347

$ws.Range("B25").Value = 30.705117562064633
$ws.Range("C25").Value = 0.1072763223178545
$ws.Range("D25").Value = 10.244718909896983

$ws.Range("B26").Value = 54.63355775188014
$ws.Range("C26").Value = 0.90719342744818043
$ws.Range("D26").Value = 7.2340881548668969

$ws.Range("B27").Value = 76.5749533337125
$ws.Range("C27").Value = 0.90265437128624515
$ws.Range("D27").Value = 9.1361400012095579

$ws.Range("B28").Value = 51.046000514598582
$ws.Range("C28").Value = 0.43378526977865395
$ws.Range("D28").Value = 6.5533922045271842
